$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.218.17"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.448.50"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "642.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.22%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.990"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("D11").Value = "3.446.69"
$ws.Range("E11").Value = "  +4.23%  "
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "94.997.46"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "4.098.48"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "3.448.95"
$ws.Range("E19").Value = "  +4.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.11%  "
$ws.Range("E22").Value = "  +7.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "501.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000191"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "91.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.627.96"
$ws.Range("E29").Value = "  +4.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.26%  "
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.59%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.566"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "512.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.150"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.82%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.99%  "
